$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "a"
$ws.Range("B3").Value = 0.0
